$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D27 was blank; fill in the missing idPagamento value. A leading apostrophe
# forces Excel to keep the all-digit value as text (matching the other
# idPagamento entries in column D), instead of auto-converting to a number.
$ws.Range("D27").Value = "'76958526368"

# Append a new registration row (row 33) mirroring the existing "Vitor Ito"
# rows above it.
$ws.Range("A33").Value = "Vitor Ito"
$ws.Range("B33").Value = 1578424633
$ws.Range("C33").Value = "'11988776655"
$ws.Range("D33").Value = "'"
$ws.Range("E33").Value = 1
$ws.Range("F33").Value = 2
$ws.Range("G33").Value = 3
$ws.Range("H33").Value = 4
$ws.Range("I33").Value = 5
$ws.Range("J33").Value = 6
$ws.Range("K33").Value = 7
$ws.Range("L33").Value = 8
$ws.Range("M33").Value = 9
$ws.Range("N33").Value = 10
$ws.Range("O33").Value = "Não"
